# Update "想去人数" (want-to-go count) values in F column across sheets
# to reflect the latest generated data (commit: Update gh-pages to output
# generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1019
$ws1.Range("F6").Value = 2203
$ws1.Range("F8").Value = 664
$ws1.Range("F10").Value = 188
$ws1.Range("F15").Value = 1301
$ws1.Range("F18").Value = 189

# --- Sheet "本地生活" (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 1983

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1983
$ws4.Range("F13").Value = 1019
$ws4.Range("F17").Value = 2203
$ws4.Range("F22").Value = 664
$ws4.Range("F24").Value = 188
$ws4.Range("F31").Value = 1301
$ws4.Range("F36").Value = 189
